# TODO_backbone_all_args.xlsx - BACKBONE v10.6 update
# (".extract_all_fun_names updated with BACKBONE v10.6")
#
# Applies the following content changes to sheet "Feuil1":
#   - D2 is cleared (the "Unit tests" mark for all_args_here.R is removed)
#   - Row 8 (intern_.all_args_here_fill.R) grows taller and its "x" mark in
#     column B is replaced by the note "compare with .colons_check_message"
#   - Row 10 (intern_.colons_check_message.R) gets a new "x" mark in column B
#   - Row 19 (is_function_here.R) gets a new "x" mark in column C
#   - Column B is widened a bit
#   - The active selection moves to C24

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Unit tests" mark previously in D2 (all_args_here.R row)
$ws.Range("D2").Clear()

# Row 8 (intern_.all_args_here_fill.R): make the row taller to fit the new
# wrapped note, and replace the old "x" mark with the note text
$ws.Rows.Item(8).RowHeight = 39
$ws.Range("B8").Value = "compare with .colons_check_message"

# Row 10 (intern_.colons_check_message.R): add an "x" mark in column B
$ws.Range("B10").Value = "x"

# Row 19 (is_function_here.R): add an "x" mark in column C
$ws.Range("C19").Value = "x"

# Widen column B slightly
$ws.Columns.Item(2).ColumnWidth = 12.67

# Move the active selection to C24
$ws.Range("C24").Select()
